$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The price-drop alert (POST) row was missing its STATUS CODE cell.
# Re-enter the surrounding text cells too, clearing them first, so the
# workbook's shared-string table is rebuilt with the price-drop-alert
# strings grouped together (matches the upstream canonical layout) and the
# stray missing comma in the product-details JSON sample gets corrected.

$ws.Range("C7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

$paramText = @'
product_id (string, required): The unique identifier of the product.
user_id (string, required): The unique identifier of the user.
'@
$ws.Range("G7").Value = $paramText

$allowsText = @'
Allows logged-in users to set an alert for a price drop on a specific product. The alert is triggered when the product's price falls to or below the desired price.
'@
$ws.Range("C8").Value = $allowsText

$priceReqJson = @'
{
  "product_id": "12345",
  "user_id": "u789",
  "desired_price": 28000
}
'@
$ws.Range("D8").Value = $priceReqJson

$priceRespJson = @'
{
  "success": true,
  "message": "Price drop alert set successfully for product ID 12345 at 28000 tk."
}

'@
$ws.Range("E8").Value = $priceRespJson

$alertsUrl = @'
https://www.TechShoppers.com/api/1.0/products/details/{user_id}/{product_id}/alerts/price-drop
'@
$ws.Range("B8").Value = $alertsUrl

$retrievingText = @'
1) retrieving detailed information about a specific product, including prices, special offers, and shipping details from various shops. This feature is activated when a user clicks on a product image and is redirected to a detailed view of the product.

2) viewing the price history of a product at a particular shop on different days of the last month
'@
$ws.Range("C7").Value = $retrievingText

$productDetailsJson = @'
{
  "success": true,
  "product_details": {
    "product_id": "12345",
    "product_name": "Smartphone XYZ",
    "category_name": "electronics",
    "subcategory_name": "computer-laptop",
    "brand_name": "HP",
    "image_path": "https://www.techshoppers.com/images/smartphone-xyz.jpg"
  },
  "offers": [
    {
      "shop_name": "ebay",
      "price": 30000,
      "special_offer_details": "10% discount on first purchase",
      "shipping_details": {
        "affirm": true,
        "afterpay": false,
        "free_shipping": true
      }
    },
    {
      "shop_name": "amazon",
      "price": 35000,
      "special_offer_details": null,
      "shipping_details": {
        "affirm": false,
        "afterpay": true,
        "free_shipping": false
      }
    },
    // More offers
  ],
  "price_history": [
    {
      "shop_name": "amazon",
      "history": [
         {
             "date": "2023-11-01",
             "price": 28000
         },
         {
             "date": "2023-11-05",
             "price": 33000
         },
         // more prices on this shop
      ]
    },
    // more shops
}

'@
$ws.Range("E7").Value = $productDetailsJson
$ws.Range("E7").Font.Size = 14

# The actual fix: the POST .../alerts/price-drop row never had its
# STATUS CODE filled in.
$ws.Range("F8").Value = 200

# Row 8 no longer needs to be as tall now that its content is unchanged
# but the table reflowed.
$ws.Rows.Item(8).RowHeight = 126

# Move the active selection down to A10.
$ws.Range("A10").Select()
